# Apply crypto price/volume updates per commit diff (GitHub Actions refresh).
# Price cells (column D) are stored as literal text in the workbook (e.g.
# "26.50", "0.0970"), so for values Excel would otherwise auto-convert to a
# number (and strip formatting like trailing zeros) we set NumberFormat to
# "@" (Text) first to force literal text entry, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.795.15'
$ws.Range("E2").Value = '  +4.74%  '
$ws.Range("D3").Value = '2.282.69'
$ws.Range("E3").Value = '  +2.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.45'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("E6").Value = '  +0.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.61'
$ws.Range("E7").Value = '  +6.98%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.427'
$ws.Range("E9").Value = '  +5.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0979'
$ws.Range("E10").Value = '  +8.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.63'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.50'
$ws.Range("E12").Value = '  +16.91%  '
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").Value = '2.618.72'
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.81'
$ws.Range("E15").Value = '  +1.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.93'
$ws.Range("E16").Value = '  +5.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.822'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").Value = '2.284.38'
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("D19").Value = '43.665.94'
$ws.Range("E19").Value = '  +4.71%  '
$ws.Range("D20").Value = '0.0₃0960'
$ws.Range("E20").Value = '  +5.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.48'
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.63'
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("E25").Value = '  +6.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.47'
$ws.Range("E26").Value = '  +6.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.01'
$ws.Range("E27").Value = '  +3.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.94'
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.139'
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("E30").Value = '  +3.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.45'
$ws.Range("E31").Value = '  +3.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.77'
$ws.Range("E32").Value = '  +6.32%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E34").Value = '  +4.81%  '
$ws.Range("E35").Value = '  +6.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.74'
$ws.Range("E36").Value = '  +1.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.83'
$ws.Range("E37").Value = '  +4.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.81'
$ws.Range("E38").Value = '  +5.66%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0248'
$ws.Range("E40").Value = '  +3.43%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.08'
$ws.Range("E42").Value = '  +28.61%  '
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.49'
$ws.Range("E44").Value = '  -0.97%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.57'
$ws.Range("E45").Value = '  +2.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.23'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0970'
$ws.Range("E47").Value = '  +0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.24'
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("D49").Value = '1.490.05'
$ws.Range("E49").Value = '  +1.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.95'
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.34'
$ws.Range("E51").Value = '  +1.72%  '
